$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (forecast) and C (hour) for rows 2-20,
# refreshed from the underlying data source.
$data = @(
    @{ Row = 2;  B = 763.0999999999999; C = 12 },
    @{ Row = 3;  B = 315;               C = 12 },
    @{ Row = 4;  B = 11;                C = 12 },
    @{ Row = 5;  B = 115;               C = 12 },
    @{ Row = 6;  B = 53;                C = 12 },
    @{ Row = 7;  B = 95;                C = 12 },
    @{ Row = 8;  B = 38;                C = 12 },
    @{ Row = 9;  B = 82;                C = 12 },
    @{ Row = 10; B = 527;               C = 12 },
    @{ Row = 11; B = 199;               C = 12 },
    @{ Row = 12; B = 738.4;             C = 12 },
    @{ Row = 13; B = 473;               C = 12 },
    @{ Row = 14; B = 727;               C = 12 },
    @{ Row = 15; B = 185;               C = 12 },
    @{ Row = 16; B = 119;               C = 12 },
    @{ Row = 17; B = 95;                C = 12 },
    @{ Row = 18; B = 4;                 C = 12 },
    @{ Row = 19; B = 11;                C = 12 },
    @{ Row = 20; B = 74;                C = 12 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
